# Apply the "break out stock.yaml completed" update to the weekly IEX.NS
# stock-history sheet:
#   - Q60 backfilled from 1 -> 0
#   - O347 backfilled from 0 -> 2
#   - R349 / R350 ("backup" column) backfilled from blank -> 0
#   - nine new weekly rows (351-359) appended with OHLCV + derived columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- backfilled corrections on existing rows -------------------------------
$ws.Cells.Item(60, 17).Value = 0   # Q60: detect_structure
$ws.Cells.Item(347, 15).Value = 2  # O347: isPivot
$ws.Cells.Item(349, 18).Value = 0  # R349: backup
$ws.Cells.Item(350, 18).Value = 0  # R350: backup

# --- newly appended weekly rows --------------------------------------------
$newRows = @(
    @(351, 45474, 181.9900054931641, 191.1999969482422, 181,               184.3500061035156, 184.3500061035156,  88221863, 2024, 7,  1, 0, 0, 0, 27, 0, 0, 0),
    @(352, 45481, 185,               185.4900054931641, 170.4700012207031, 177.1300048828125, 177.1300048828125,  65255091, 2024, 7,  8, 0, 0, 0, 28, 0, 0, 0),
    @(353, 45488, 178.25,            182.3699951171875, 168.5,             169.0700073242188, 169.0700073242188,  38097836, 2024, 7, 15, 0, 0, 0, 29, 0, 1, 1),
    @(354, 45495, 168,               178.6399993896484, 159.3500061035156, 176.6600036621094, 176.6600036621094,  54937155, 2024, 7, 22, 0, 0, 0, 30, 2, 0, 0),
    @(355, 45502, 178.5,             197.8000030517578, 178.0500030517578, 195.3200073242188, 195.3200073242188, 120973318, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(356, 45509, 192.8899993896484, 200.8999938964844, 185.3200073242188, 192.6699981689453, 192.6699981689453,  92516865, 2024, 8,  5, 0, 0, 0, 32, 0, 0, 0),
    @(357, 45516, 190,               195.6199951171875, 184.4100036621094, 194.8200073242188, 194.8200073242188,  38946437, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
    @(358, 45523, 196.3999938964844, 200.6900024414062, 188.25,            188.9700012207031, 188.9700012207031,  47448409, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
    @(359, 45530, 190.7599945068359, 208.7599945068359, 186.3000030517578, 203.6300048828125, 203.6300048828125, 127066520, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value  = $row[1]   # A Datetime (serial)
    $ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item(350, 1).NumberFormat  # match existing date style
    $ws.Cells.Item($r, 2).Value  = $row[2]   # B Open
    $ws.Cells.Item($r, 3).Value  = $row[3]   # C High
    $ws.Cells.Item($r, 4).Value  = $row[4]   # D Low
    $ws.Cells.Item($r, 5).Value  = $row[5]   # E Close
    $ws.Cells.Item($r, 6).Value  = $row[6]   # F Adj Close
    $ws.Cells.Item($r, 7).Value  = $row[7]   # G Volume
    $ws.Cells.Item($r, 8).Value  = $row[8]   # H Year
    $ws.Cells.Item($r, 9).Value  = $row[9]   # I Month
    $ws.Cells.Item($r, 10).Value = $row[10]  # J Day
    $ws.Cells.Item($r, 11).Value = $row[11]  # K Hour
    $ws.Cells.Item($r, 12).Value = $row[12]  # L Minute
    $ws.Cells.Item($r, 13).Value = $row[13]  # M Second
    $ws.Cells.Item($r, 14).Value = $row[14]  # N Week
    $ws.Cells.Item($r, 15).Value = $row[15]  # O isPivot
    $ws.Cells.Item($r, 16).Value = $row[16]  # P two_line_structure
    $ws.Cells.Item($r, 17).Value = $row[17]  # Q detect_structure
    # R (backup) intentionally left blank for the newest, not-yet-backed-up rows
}
